# standardized_IS.xlsx - Published fix
# 1. fixed IS: negative Interest Expense_12m  (no visible cell-data change in this diff)
# 2. validated IS: Net Income (loss) (operations)_12m
#      -> relabel the "Income Before Income Tax / IncomeTaxExpenseBenefit" row to
#         "Provision for (benefit) Income Tax"
#      -> relabel the "Net Income (loss) (operations) / ProfitLoss" row to
#         "Net Income (loss) (continous operations)"
# 3. removed IS: Net Income (loss) to parent (incl. Non contr. Interest)_12m
#      -> delete the NetIncomeLoss row entirely
#    also removes the now-orphaned Domestic/Foreign "Income Before Income Tax" rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Domestic/Foreign "Income (Loss) from Continuing Operations
# before Income Taxes" detail rows (old rows 18 & 19).
$ws.Range("A18:A19").EntireRow.Delete()

# After the delete above, the remaining rows shift up by two:
#   old row 21 (IncomeTaxExpenseBenefit)            -> now row 19
#   old row 22 (ProfitLoss)                          -> now row 20
#   old row 23 (NetIncomeLoss, to be removed)         -> now row 21
$ws.Range("A19").Value = "Provision for (benefit) Income Tax"
$ws.Range("A20").Value = "Net Income (loss) (continous operations)"

# Drop the "Net Income (loss) to parent (incl. Non contr. Interest)" row.
$ws.Range("A21:E21").EntireRow.Delete()

# Match the author's final selection (the now-blank trailing row).
$ws.Rows("21:21").Select()
